$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a plain-numeric-looking price string must be pre-formatted as
# Text so Excel keeps them as strings (matching the workbook's inlineStr cells)
# instead of silently converting them to numbers.
$textCells = @('D5','D8','D9','D10','D11','D13','D14','D16','D17','D18','D19','D21','D22','D24','D25','D26','D27','D29','D30','D32','D33','D34','D35','D36','D37','D39','D40','D41','D42','D43','D44','D45','D46','D47','D48','D49','D51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.772.50'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '1.784.01'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '310.48'
$ws.Range('E5').Value = '  -2.05%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D8').Value = '0.3852'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').Value = '0.07814'
$ws.Range('E9').Value = '  -7.99%  '
$ws.Range('D10').Value = '1.088'
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('D11').Value = '40.76'
$ws.Range('E11').Value = '  -2.76%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '6.193'
$ws.Range('E13').Value = '  -3.77%  '
$ws.Range('D14').Value = '20.14'
$ws.Range('E14').Value = '  -4.28%  '
$ws.Range('D15').Value = '1.780.01'
$ws.Range('E15').Value = '  -2.08%  '
$ws.Range('D16').Value = '7.204'
$ws.Range('E16').Value = '  -4.42%  '
$ws.Range('D17').Value = '91.32'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '0.00001074'
$ws.Range('E18').Value = '  -5.79%  '
$ws.Range('D19').Value = '0.06552'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '17.00'
$ws.Range('E21').Value = '  -4.31%  '
$ws.Range('D22').Value = '5.901'
$ws.Range('E22').Value = '  -3.13%  '
$ws.Range('D23').Value = '27.819.18'
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('D24').Value = '10.99'
$ws.Range('E24').Value = '  -4.16%  '
$ws.Range('D25').Value = '2.226'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '159.81'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').Value = '20.17'
$ws.Range('E27').Value = '  -4.14%  '
$ws.Range('D28').Value = '1.984.24'
$ws.Range('E28').Value = '  -2.16%  '
$ws.Range('D29').Value = '2.363'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').Value = '123.23'
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').Value = '1.031'
$ws.Range('E32').Value = '  -6.23%  '
$ws.Range('D33').Value = '3.634'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').Value = '5.477'
$ws.Range('E34').Value = '  -4.38%  '
$ws.Range('D35').Value = '0.07051'
$ws.Range('E35').Value = '  -5.36%  '
$ws.Range('D36').Value = '0.02300'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('D37').Value = '8.785'
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('E38').Value = '  -5.05%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '4.991'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '11.43'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').Value = '0.6075'
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '1.151'
$ws.Range('E43').Value = '  -3.74%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '13.08'
$ws.Range('E44').Value = '  -3.99%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '1.317'
$ws.Range('E45').Value = '  -5.93%  '
$ws.Range('D46').Value = '0.5892'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').Value = '3.701'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').Value = '125.81'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').Value = '1.199'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('E50').Value = '  -4.72%  '
$ws.Range('D51').Value = '0.06840'
$ws.Range('E51').Value = '  -1.96%  '
